# Applies the "Game Set Expenses" budget update to the Rookie Team Season Budget
# worksheet (Sheet1) of the FTC 2023-2024 budget workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Game Supplies section: add a new line item in the previously-blank row 14
#    for servo / wiring supplies, with Budget and Actual amounts of $313.92.
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Connectors, Wires, Locktight, Surgical tubing, Compliant wheels, Servo motor, Servo programmer, Servo block"
$ws.Range("B14").Value = 313.92
$ws.Range("C14").Value = 313.92

# ---------------------------------------------------------------------------
# 2) Row 16 ("Scoring Elements" label) is cleared out - the line item stays
#    blank (category/rationale text remain, only the item name is removed).
# ---------------------------------------------------------------------------
$ws.Range("A16").ClearContents()

# ---------------------------------------------------------------------------
# 3) Grants section: rename the "IBM " donor line to clarify the grant terms.
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "IBM grant only for purchases from FIRST"

# ---------------------------------------------------------------------------
# 4) Optional Expenses subtotal (row 27): the Budget Amt. column formula now
#    starts at row 22 instead of row 21 (the IBM grant row is excluded from
#    the optional-expenses budget total).
# ---------------------------------------------------------------------------
$ws.Range("B27").Formula = "=SUM(B22:B26)"

# ---------------------------------------------------------------------------
# 5) Bottom Line section: insert a new row 40 ("Available Funds for Purchase
#    from FIRST only") above the existing "Potential Optional Expenses" /
#    "Credit Deficit" rows, pushing them down to rows 41 and 42.
# ---------------------------------------------------------------------------
$ws.Range("A40").EntireRow.Insert()

# Copy formatting from the row above (row 39) onto the newly inserted row 40
# so it keeps the same borders/number formats as the rest of the section.
$ws.Range("A39:E39").Copy()
$ws.Range("A40:E40").PasteSpecial(-4122)
$ws.Rows.Item(40).RowHeight = $ws.Rows.Item(39).RowHeight
$excel.CutCopyMode = 0

# New row 40 content: a hard-coded available-funds figure (not a formula).
$ws.Range("A40").Value = "Available Funds for Purchase from FIRST only"
$ws.Range("B40").Value = 650
$ws.Range("C40").ClearContents()

# Former row 40 ("Potential Optional Expenses"), now row 41, is relabeled.
$ws.Range("A41").Value = "Game Set Expenses"

# Former row 41 ("Credit/Deficit"), now row 42, keeps its label/formulas but
# the Credit/Deficit formulas must reference the shifted "Game Set Expenses"
# row (41) instead of the old row 40.
$ws.Range("A42").Value = "Credit/Deficit"
$ws.Range("B42").Formula = "=B39-(B38+B41)"
$ws.Range("C42").Formula = "=C39-(C38+C41)"

# ---------------------------------------------------------------------------
# 6) NOTE: the footer note ("Updated 10.22.2023") automatically moved from
#    row 44 down to row 45 as a side-effect of the row 40 insertion above
#    (everything below row 40 shifts down by one row), so no further action
#    is required here - row 44 is now blank and the note lives in row 45.
# ---------------------------------------------------------------------------

$wb.Save()
